$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cell E8 from "Good Morning" to "GIT UPDATE"
$ws.Range("E8").Value = "GIT UPDATE"

# Reflect the selection change recorded in the saved file (user selected E8)
$ws.Range("E8").Select()
